$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Jos Buttler <dagger>" uses a NBSP (U+00A0) before the dagger glyph, matching
# the byte sequence already used for this label in the existing rows (2-8).
$batsman = "Jos Buttler" + [char]0x00A0 + [char]0x2020

# New match rows to append (rows 9-15), columns A-K:
# venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$rows = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Rajasthan Royals", "Delhi Capitals", $batsman, "22", "9", "3", "1", "244.44"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Rajasthan Royals", "Kolkata Knight Riders", $batsman, "21", "16", "1", "2", "131.25"),
    @(" Sharjah", " September 27 2020", "Royals won by 4 wickets (with 3 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", $batsman, "4", "7", "0", "0", "57.14"),
    @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", $batsman, "16", "13", "1", "1", "123.07"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Rajasthan Royals", "Mumbai Indians", $batsman, "70", "44", "4", "5", "159.09"),
    @(" Abu Dhabi", " October 03 2020", "RCB won by 8 wickets (with 5 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", $batsman, "22", "12", "3", "1", "183.33"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Rajasthan Royals", "Delhi Capitals", $batsman, "13", "8", "2", "0", "162.50")
)

$startRow = 9
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text storage so values like "13" stay as strings, not numbers
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
    }
}

Write-Host "done"
